$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8986206666666666
$ws.Range("H2").Value = 2.695862
$ws.Range("I2").Value = 0.1661804693926261
$ws.Range("J2").Value = 0.1661804693926262
$ws.Range("O2").Value = 0.858667536176972
$ws.Range("P2").Value = 0.858667536176972
$ws.Range("Q2").Value = 0.1411975694913333
$ws.Range("R2").Value = 1.270778125422
$ws.Range("S2").Value = 0.142693774214099
$ws.Range("T2").Value = 0.142693774214099

$ws.Range("G3").Value = 0.8986206666666666
$ws.Range("H3").Value = 2.695862
$ws.Range("I3").Value = 0.1661804693926261
$ws.Range("J3").Value = 0.1661804693926262
$ws.Range("Q3").Value = 0.02324042722155556
$ws.Range("R3").Value = 0.209163844994
$ws.Range("S3").Value = 0.02348669517852713
$ws.Range("T3").Value = 0.02348669517852714

$ws.Range("I4").Value = 0.3901029163453022
$ws.Range("J4").Value = 0.3901029163453023
$ws.Range("O4").Value = 0.858667536176972
$ws.Range("P4").Value = 0.858667536176972
$ws.Range("Q4").Value = 0.331456421087
$ws.Range("S4").Value = 0.3349687100336721
$ws.Range("T4").Value = 0.3349687100336721

$ws.Range("I5").Value = 0.3901029163453022
$ws.Range("J5").Value = 0.3901029163453023
$ws.Range("S5").Value = 0.05513420631163012
$ws.Range("T5").Value = 0.05513420631163013

$ws.Range("G6").Value = 2.399397
$ws.Range("H6").Value = 7.198191
$ws.Range("I6").Value = 0.4437166142620716
$ws.Range("J6").Value = 0.4437166142620716
$ws.Range("O6").Value = 0.858667536176972
$ws.Range("P6").Value = 0.858667536176972
$ws.Range("Q6").Value = 0.377010052419
$ws.Range("R6").Value = 3.393090471771
$ws.Range("S6").Value = 0.3810050519292009
$ws.Range("T6").Value = 0.3810050519292009

$ws.Range("G7").Value = 2.399397
$ws.Range("H7").Value = 7.198191
$ws.Range("I7").Value = 0.4437166142620716
$ws.Range("J7").Value = 0.4437166142620716
$ws.Range("Q7").Value = 0.062054005013
$ws.Range("R7").Value = 0.558486045117
$ws.Range("S7").Value = 0.06271156233287067
$ws.Range("T7").Value = 0.06271156233287067

